# Applies updated crypto price/volume data to sheet1 (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.171.16"
$ws.Range("E2").Value = "  +4.17%  "
$ws.Range("D3").Value = "2.489.02"
$ws.Range("E3").Value = "  +2.80%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.82"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.13"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.97%  "
$ws.Range("E7").Value = "  +2.25%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +2.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.17"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0815"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.52"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("D15").Value = "2.878.57"
$ws.Range("E15").Value = "  +2.86%  "
$ws.Range("D16").Value = "2.493.81"
$ws.Range("E16").Value = "  +3.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.856"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("D18").Value = "47.094.70"
$ws.Range("E18").Value = "  +4.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.92"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.41%  "
$ws.Range("E20").Value = "  +5.19%  "
$ws.Range("D21").Value = "0.0₃0940"
$ws.Range("E21").Value = "  +2.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.76"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.48"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "250.55"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.52%  "
$ws.Range("E25").Value = "  +3.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.22"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.80%  "
$ws.Range("E29").Value = "  +0.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.97"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.18%  "
$ws.Range("E31").Value = "  +8.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.42"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.74"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.30%  "
$ws.Range("E34").Value = "  +4.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0793"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.20%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.99"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.72"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.98"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.60%  "
$ws.Range("E40").Value = "  +2.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "122.56"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.55%  "
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("E43").Value = "  +3.13%  "
$ws.Range("E44").Value = "  +3.09%  "
$ws.Range("D45").Value = "1.969.96"
$ws.Range("E45").Value = "  +1.41%  "
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("E50").Value = "  +9.54%  "
$ws.Range("E51").Value = "  +3.23%  "
